# Generate Report for Handoff
#
# Refresh the "Latest Handoff" / "Latest HO Xliff Generate Date" timestamps
# for the 3df4e011-dc30-43c6-aa5c-d3521ff6d82b.md file after a new handoff
# report was generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 3df4e011-... row (row 7 of the table / row 5 of the worksheet).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-08-13 20:54:22"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the same file.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-08-13 20:54:14"

# de-de sheet: "Latest Handoff Datetime" column (H) for the same file.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-08-13 20:54:22"
